$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footer: merge dash runs into one ---
$f = $sec.Footers(1)
$findDashes = "---------------------------------------------------------------------------------------------------------" + "-------------------------"
$found1 = $f.Range.Find.Execute($findDashes, $true, $false, $false, $false, $false, $true, 1, $false, $findDashes, 2)
Write-Output "dash merge found=$found1"

# --- Footer: merge "01 BP..." runs into one (keep leading 3-space run separate) ---
$findAddr = "01 BP: 613,   Porto " + [char]0x2013 + " Novo, B" + [char]0x00E9 + "nin   e-mail: secretariat@imsp-uac.org   site web: www.imsp-benin.com "
$found2 = $f.Range.Find.Execute($findAddr, $true, $false, $false, $false, $false, $true, 1, $false, $findAddr, 2)
Write-Output "addr merge found=$found2"

# --- Header: merge asterisk runs into one ---
$h = $sec.Headers(1)
$findStars = "                                                ****************"
$found3 = $h.Range.Find.Execute($findStars, $true, $false, $false, $false, $false, $true, 1, $false, $findStars, 2)
Write-Output "stars merge found=$found3"
